# Fruta / hortaliza, semanal
# Insert the new weekly price record as row 63 (pushing the existing
# rows 63-65 down to 64-66) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 63:65 down to 64:66 and leave a blank row 63 to populate.
$ws.Rows(63).Insert()

# Populate the newly inserted row 63 with the new record.
$ws.Cells.Item(63, 1).Value = 9
$ws.Cells.Item(63, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(63, 3).Value = "Metropolitana"
$ws.Cells.Item(63, 4).Value = 44568
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = "Fruta"
$ws.Cells.Item(63, 7).Value = 100101
$ws.Cells.Item(63, 8).Value = "Berries"
$ws.Cells.Item(63, 9).Value = 100101004
$ws.Cells.Item(63, 10).Value = "Frambuesa"
$ws.Cells.Item(63, 11).Value = "Sin especificar"
$ws.Cells.Item(63, 12).Value = "Primera"
$ws.Cells.Item(63, 13).Value = 680
$ws.Cells.Item(63, 14).Value = 7500
$ws.Cells.Item(63, 15).Value = 8000
$ws.Cells.Item(63, 16).Value = 7757
$ws.Cells.Item(63, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(63, 18).Value = "Provincia de Linares"
$ws.Cells.Item(63, 19).Value = 3878
$ws.Cells.Item(63, 20).Value = 2
